# Insert a new weekly price record at row 241 ("Fruta / hortaliza, semanal").
# This pushes the existing rows 241-276 down to 242-277 (dimension grows
# from A1:R276 to A1:R277) and fills the newly inserted row 241 with the
# latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 241, shifting rows 241:276 -> 242:277
$ws.Rows.Item(241).Insert()

# Populate the new row 241 with the latest record
$ws.Cells.Item(241, 1).Value = 7
$ws.Cells.Item(241, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(241, 3).Value = "Ñuble"
$ws.Cells.Item(241, 4).Value = 44505
$ws.Cells.Item(241, 5).Value = 16
$ws.Cells.Item(241, 6).Value = 100114014
$ws.Cells.Item(241, 7).Value = "Betarraga"
$ws.Cells.Item(241, 8).Value = "Sin especificar"
$ws.Cells.Item(241, 9).Value = "Primera"
$ws.Cells.Item(241, 10).Value = 400
$ws.Cells.Item(241, 11).Value = 700
$ws.Cells.Item(241, 12).Value = 800
$ws.Cells.Item(241, 13).Value = 750
$ws.Cells.Item(241, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(241, 15).Value = "Región del Maule"
$ws.Cells.Item(241, 16).Value = 150
$ws.Cells.Item(241, 17).Value = 5
$ws.Cells.Item(241, 18).Value = "Hortaliza"
